# Re-experiment of 5-fold grouping: update computed accuracy values (and the
# two cells whose number format changed alongside their values), move the
# active selection, and nudge the saved window position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated 5-fold test-accuracy figures -----------------------------
$ws.Range("B3").Value = 0.72897199999999995
$ws.Range("B4").Value = 0.73831800000000003

$ws.Range("B5").Value = 0.77570099999999997
$ws.Range("B5").NumberFormat = "0.0000%"

$ws.Range("B6").Value = 0.71028000000000002
$ws.Range("B6").NumberFormat = "0.000%"

$ws.Range("C6").Value = 0.97663599999999995
$ws.Range("C6").NumberFormat = "0.0000%"

$ws.Range("B7").Value = 0.76635500000000001

# --- Selection moves from the old C17 to B3 ----------------------------
[void]$ws.Range("B3").Select()

# --- Saved window position (best effort - host may not persist this) --
$win = $wb.Windows.Item(1)
$win.Left = 7360
$win.Top = 900
$excel.Left = 7360
$excel.Top = 900
